$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 193
$ws.Cells.Item(193, 1).Value = 192.0
$ws.Cells.Item(193, 2).Value = 'Monday, Jan 09'
$ws.Cells.Item(193, 3).Value = '2:20 PM'
$ws.Cells.Item(193, 4).Value = 'LO3905'
$ws.Cells.Item(193, 5).Value = 'Warsaw'
$ws.Cells.Item(193, 6).Value = '(WAW)'
$ws.Cells.Item(193, 7).Value = 'LOT '
$ws.Cells.Item(193, 8).Value = 'E190'
$ws.Cells.Item(193, 9).Value = '(SP-LMH)'
$ws.Cells.Item(193, 10).Value = '2:39 PM'
$ws.Cells.Item(193, 12).Value = '0 hours, 19 minutes'

# Row 194
$ws.Cells.Item(194, 1).Value = 193.0
$ws.Cells.Item(194, 2).Value = 'Monday, Jan 09'
$ws.Cells.Item(194, 3).Value = '2:25 PM'
$ws.Cells.Item(194, 4).Value = 'FR4024'
$ws.Cells.Item(194, 5).Value = 'Lisbon'
$ws.Cells.Item(194, 6).Value = '(LIS)'
$ws.Cells.Item(194, 7).Value = 'Buzz '
$ws.Cells.Item(194, 8).Value = 'B38M'
$ws.Cells.Item(194, 9).Value = '(SP-RZD)'
$ws.Cells.Item(194, 10).Value = '2:20 PM'
$ws.Cells.Item(194, 12).Value = '0 hours, -5 minutes'

# Row 195
$ws.Cells.Item(195, 1).Value = 194.0
$ws.Cells.Item(195, 2).Value = 'Monday, Jan 09'
$ws.Cells.Item(195, 3).Value = '2:45 PM'
$ws.Cells.Item(195, 4).Value = 'FR7954'
$ws.Cells.Item(195, 5).Value = 'Prague'
$ws.Cells.Item(195, 6).Value = '(PRG)'
$ws.Cells.Item(195, 7).Value = 'Ryanair '
$ws.Cells.Item(195, 8).Value = 'B738'
$ws.Cells.Item(195, 9).Value = '(SP-RKC)'
$ws.Cells.Item(195, 10).Value = '2:36 PM'
$ws.Cells.Item(195, 12).Value = '0 hours, -9 minutes'

# Row 196
$ws.Cells.Item(196, 1).Value = 195.0
$ws.Cells.Item(196, 2).Value = 'Monday, Jan 09'
$ws.Cells.Item(196, 3).Value = '2:50 PM'
$ws.Cells.Item(196, 4).Value = 'FR6364'
$ws.Cells.Item(196, 5).Value = 'Shannon'
$ws.Cells.Item(196, 6).Value = '(SNN)'
$ws.Cells.Item(196, 7).Value = 'Ryanair '
$ws.Cells.Item(196, 8).Value = 'B738'
$ws.Cells.Item(196, 9).Value = '(EI-EBP)'
$ws.Cells.Item(196, 10).Value = '2:32 PM'
$ws.Cells.Item(196, 12).Value = '0 hours, -18 minutes'

# Row 197
$ws.Cells.Item(197, 1).Value = 196.0
$ws.Cells.Item(197, 2).Value = 'Monday, Jan 09'
$ws.Cells.Item(197, 3).Value = '3:00 PM'
$ws.Cells.Item(197, 4).Value = 'LG5741'
$ws.Cells.Item(197, 5).Value = 'Luxembourg'
$ws.Cells.Item(197, 6).Value = '(LUX)'
$ws.Cells.Item(197, 7).Value = 'Luxair '
$ws.Cells.Item(197, 8).Value = 'DH8D'
$ws.Cells.Item(197, 9).Value = '(LX-LGE)'
$ws.Cells.Item(197, 10).Value = '2:41 PM'
$ws.Cells.Item(197, 12).Value = '0 hours, -19 minutes'

# Row 198
$ws.Cells.Item(198, 1).Value = 197.0
$ws.Cells.Item(198, 2).Value = 'Monday, Jan 09'
$ws.Cells.Item(198, 3).Value = '3:01 PM'
$ws.Cells.Item(198, 4).Value = 'LPR41'
$ws.Cells.Item(198, 5).Value = 'Szczecin'
$ws.Cells.Item(198, 6).Value = '(SZZ)'
$ws.Cells.Item(198, 7).Value = 'Polish Medical Air Rescue '
$ws.Cells.Item(198, 8).Value = 'LJ75'
$ws.Cells.Item(198, 9).Value = '(SP-MXR)'
$ws.Cells.Item(198, 10).Value = '4:22 PM'
$ws.Cells.Item(198, 12).Value = '1 hours, 21 minutes'

# Row 199
$ws.Cells.Item(199, 1).Value = 198.0
$ws.Cells.Item(199, 2).Value = 'Monday, Jan 09'
$ws.Cells.Item(199, 3).Value = '3:15 PM'
$ws.Cells.Item(199, 4).Value = 'FR3721'
$ws.Cells.Item(199, 5).Value = 'Billund'
$ws.Cells.Item(199, 6).Value = '(BLL)'
$ws.Cells.Item(199, 7).Value = 'Ryanair '
$ws.Cells.Item(199, 8).Value = 'B738'
$ws.Cells.Item(199, 9).Value = '(9H-QBX)'
$ws.Cells.Item(199, 10).Value = '3:07 PM'
$ws.Cells.Item(199, 12).Value = '0 hours, -8 minutes'

# Row 200
$ws.Cells.Item(200, 1).Value = 199.0
$ws.Cells.Item(200, 2).Value = 'Monday, Jan 09'
$ws.Cells.Item(200, 3).Value = '3:20 PM'
$ws.Cells.Item(200, 4).Value = 'U23815'
$ws.Cells.Item(200, 5).Value = 'Paris'
$ws.Cells.Item(200, 6).Value = '(CDG)'
$ws.Cells.Item(200, 7).Value = 'easyJet '
$ws.Cells.Item(200, 8).Value = 'A320'
$ws.Cells.Item(200, 9).Value = '(OE-IVS)'
$ws.Cells.Item(200, 10).Value = '3:10 PM'
$ws.Cells.Item(200, 12).Value = '0 hours, -10 minutes'

# Row 201
$ws.Cells.Item(201, 1).Value = 200.0
$ws.Cells.Item(201, 2).Value = 'Monday, Jan 09'
$ws.Cells.Item(201, 3).Value = '3:25 PM'
$ws.Cells.Item(201, 4).Value = 'FR5623'
$ws.Cells.Item(201, 5).Value = 'Edinburgh'
$ws.Cells.Item(201, 6).Value = '(EDI)'
$ws.Cells.Item(201, 7).Value = 'Buzz '
$ws.Cells.Item(201, 8).Value = 'B38M'
$ws.Cells.Item(201, 9).Value = '(SP-RZA)'
$ws.Cells.Item(201, 10).Value = '3:23 PM'
$ws.Cells.Item(201, 12).Value = '0 hours, -2 minutes'

# Row 202
$ws.Cells.Item(202, 1).Value = 201.0
$ws.Cells.Item(202, 2).Value = 'Monday, Jan 09'
$ws.Cells.Item(202, 3).Value = '3:50 PM'
$ws.Cells.Item(202, 4).Value = 'FR2332'
$ws.Cells.Item(202, 5).Value = 'Leeds'
$ws.Cells.Item(202, 6).Value = '(LBA)'
$ws.Cells.Item(202, 7).Value = 'Buzz '
$ws.Cells.Item(202, 8).Value = 'B38M'
$ws.Cells.Item(202, 9).Value = '(SP-RZF)'
$ws.Cells.Item(202, 10).Value = '3:51 PM'
$ws.Cells.Item(202, 12).Value = '0 hours, 1 minutes'

# Row 203
$ws.Cells.Item(203, 1).Value = 202.0
$ws.Cells.Item(203, 2).Value = 'Monday, Jan 09'
$ws.Cells.Item(203, 3).Value = '4:10 PM'
$ws.Cells.Item(203, 4).Value = 'KL1995'
$ws.Cells.Item(203, 5).Value = 'Amsterdam'
$ws.Cells.Item(203, 6).Value = '(AMS)'
$ws.Cells.Item(203, 7).Value = 'KLM '
$ws.Cells.Item(203, 8).Value = 'E190'
$ws.Cells.Item(203, 9).Value = '(PH-EXC)'
$ws.Cells.Item(203, 10).Value = '4:09 PM'
$ws.Cells.Item(203, 12).Value = '0 hours, -1 minutes'

# Row 204
$ws.Cells.Item(204, 1).Value = 203.0
$ws.Cells.Item(204, 2).Value = 'Monday, Jan 09'
$ws.Cells.Item(204, 3).Value = '4:15 PM'
$ws.Cells.Item(204, 4).Value = 'FR2712'
$ws.Cells.Item(204, 5).Value = 'London'
$ws.Cells.Item(204, 6).Value = '(STN)'
$ws.Cells.Item(204, 7).Value = 'Ryanair '
$ws.Cells.Item(204, 8).Value = 'B738'
$ws.Cells.Item(204, 9).Value = '(EI-DWY)'
$ws.Cells.Item(204, 10).Value = '4:06 PM'
$ws.Cells.Item(204, 12).Value = '0 hours, -9 minutes'

# Row 205
$ws.Cells.Item(205, 1).Value = 204.0
$ws.Cells.Item(205, 2).Value = 'Monday, Jan 09'
$ws.Cells.Item(205, 3).Value = '4:20 PM'
$ws.Cells.Item(205, 4).Value = 'LO3923'
$ws.Cells.Item(205, 5).Value = 'Warsaw'
$ws.Cells.Item(205, 6).Value = '(WAW)'
$ws.Cells.Item(205, 7).Value = 'LOT '
$ws.Cells.Item(205, 8).Value = 'E190'
$ws.Cells.Item(205, 9).Value = '(SP-LMF)'
$ws.Cells.Item(205, 10).Value = '4:12 PM'
$ws.Cells.Item(205, 12).Value = '0 hours, -8 minutes'

# Row 206
$ws.Cells.Item(206, 1).Value = 205.0
$ws.Cells.Item(206, 2).Value = 'Monday, Jan 09'
$ws.Cells.Item(206, 3).Value = '5:05 PM'
$ws.Cells.Item(206, 4).Value = 'W65072'
$ws.Cells.Item(206, 5).Value = 'Nice'
$ws.Cells.Item(206, 6).Value = '(NCE)'
$ws.Cells.Item(206, 7).Value = 'Wizz Air '
$ws.Cells.Item(206, 8).Value = 'A21N'
$ws.Cells.Item(206, 9).Value = '(HA-LVH)'
$ws.Cells.Item(206, 10).Value = '5:07 PM'
$ws.Cells.Item(206, 12).Value = '0 hours, 2 minutes'

# Row 207
$ws.Cells.Item(207, 1).Value = 206.0
$ws.Cells.Item(207, 2).Value = 'Monday, Jan 09'
$ws.Cells.Item(207, 3).Value = '5:10 PM'
$ws.Cells.Item(207, 4).Value = 'FR6235'
$ws.Cells.Item(207, 5).Value = 'Copenhagen'
$ws.Cells.Item(207, 6).Value = '(CPH)'
$ws.Cells.Item(207, 7).Value = 'Ryanair '
$ws.Cells.Item(207, 8).Value = 'B738'
$ws.Cells.Item(207, 9).Value = '(SP-RSA)'
$ws.Cells.Item(207, 10).Value = '4:59 PM'
$ws.Cells.Item(207, 12).Value = '0 hours, -11 minutes'

# Row 208
$ws.Cells.Item(208, 1).Value = 207.0
$ws.Cells.Item(208, 2).Value = 'Monday, Jan 09'
$ws.Cells.Item(208, 3).Value = '5:25 PM'
$ws.Cells.Item(208, 4).Value = 'U26939'
$ws.Cells.Item(208, 5).Value = 'Edinburgh'
$ws.Cells.Item(208, 6).Value = '(EDI)'
$ws.Cells.Item(208, 7).Value = 'easyJet (Europcar Livery) '
$ws.Cells.Item(208, 8).Value = 'A20N'
$ws.Cells.Item(208, 9).Value = '(G-UZHO)'
$ws.Cells.Item(208, 10).Value = '5:03 PM'
$ws.Cells.Item(208, 12).Value = '0 hours, -22 minutes'

# Row 209
$ws.Cells.Item(209, 1).Value = 208.0
$ws.Cells.Item(209, 2).Value = 'Monday, Jan 09'
$ws.Cells.Item(209, 3).Value = '6:10 PM'
$ws.Cells.Item(209, 4).Value = 'W65010'
$ws.Cells.Item(209, 5).Value = 'London'
$ws.Cells.Item(209, 6).Value = '(LGW)'
$ws.Cells.Item(209, 7).Value = 'Wizz Air '
$ws.Cells.Item(209, 8).Value = 'A21N'
$ws.Cells.Item(209, 9).Value = '(HA-LVO)'
$ws.Cells.Item(209, 10).Value = '5:40 PM'
$ws.Cells.Item(209, 12).Value = '0 hours, -30 minutes'
